$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '(''PL82 0606 5923 2589 0099 5855 6419'', ''next generation'')'
$ws.Cells.Item(2, 2).Value = '[2017.0, 30.4334, 22.8041, 18.4]'
$ws.Cells.Item(2, 3).Value = '[2017.0, 30.4334, 22.8041, 18.4]'
$ws.Cells.Item(2, 4).Value = 'PASS'
$ws.Cells.Item(2, 5).Value = '[]'
$ws.Cells.Item(3, 1).Value = '(''NO86 8247 5299 293'', '''')'
$ws.Cells.Item(3, 2).Value = '[2015.0, 4.5325, 23.36493, 20.12]'
$ws.Cells.Item(3, 3).Value = '[2015.0, 4.5325, 23.36493, 20.12]'
$ws.Cells.Item(3, 4).Value = 'PASS'
$ws.Cells.Item(3, 5).Value = '[]'
$ws.Cells.Item(4, 1).Value = '(''LV95 ROBP YYK3 7KVF 97PQ V'', '''')'
$ws.Cells.Item(4, 2).Value = '[2016.0, 21.66792, 22.75278, 25.6]'
$ws.Cells.Item(4, 3).Value = '[2016.0, 21.66792, 22.75278, 25.6]'
$ws.Cells.Item(4, 4).Value = 'PASS'
$ws.Cells.Item(4, 5).Value = '[]'
$ws.Cells.Item(5, 1).Value = '(''TR07 5268 2YWM TXNS DDJB TDYH CG'', '''')'
$ws.Cells.Item(5, 2).Value = '[2019.0, 7.53757, 8.57464, 21.1]'
$ws.Cells.Item(5, 3).Value = '[2019.0, 7.53757, 8.57464, 21.1]'
$ws.Cells.Item(5, 4).Value = 'PASS'
$ws.Cells.Item(5, 5).Value = '[]'
$ws.Cells.Item(6, 1).Value = '(''FO22 2384 9855 6956 12'', '''')'
$ws.Cells.Item(6, 2).Value = '[2006.0, 33.33669, 34.05317, 23.49]'
$ws.Cells.Item(6, 3).Value = '[2006.0, 33.33669, 34.05317, 23.49]'
$ws.Cells.Item(6, 4).Value = 'PASS'
$ws.Cells.Item(6, 5).Value = '[]'
$ws.Cells.Item(7, 1).Value = '(''AZ93 DCIN 1WE0 3N4S VCJR OPBN V9RT'', '''', ''optimizing'')'
$ws.Cells.Item(7, 2).Value = '[2011.0, 17.05823, 18.66]'
$ws.Cells.Item(7, 3).Value = '[2011.0, 17.05823, 18.66]'
$ws.Cells.Item(7, 4).Value = 'PASS'
$ws.Cells.Item(7, 5).Value = '[]'
$ws.Cells.Item(8, 1).Value = '(''LV14 EUQC Z3OS DRE2 BV3X E'', '''')'
$ws.Cells.Item(8, 2).Value = '[2017.0, 9.13007, 5.43523, 17.06]'
$ws.Cells.Item(8, 3).Value = '[2017.0, 9.13007, 5.43523, 17.06]'
$ws.Cells.Item(8, 4).Value = 'PASS'
$ws.Cells.Item(8, 5).Value = '[]'
$ws.Cells.Item(9, 1).Value = '(''MD10 QSDS C0PP T8YB NWEB MUCW'', '''')'
$ws.Cells.Item(9, 2).Value = '[2006.0, 20.3756, 22.8325, 4.83]'
$ws.Cells.Item(9, 3).Value = '[2006.0, 20.3756, 22.8325, 4.83]'
$ws.Cells.Item(9, 4).Value = 'PASS'
$ws.Cells.Item(9, 5).Value = '[]'
$ws.Cells.Item(10, 1).Value = '(''ME21 6004 3035 6245 8080 85'', '''')'
$ws.Cells.Item(10, 2).Value = '[2009.0, 29.86816, 16.76982, 18.74]'
$ws.Cells.Item(10, 3).Value = '[2009.0, 29.86816, 16.76982, 18.74]'
$ws.Cells.Item(10, 4).Value = 'PASS'
$ws.Cells.Item(10, 5).Value = '[]'
$ws.Cells.Item(11, 1).Value = '(''SM88 C084 1092 877W KBGO GYDP RB9'', ''Persevering'')'
$ws.Cells.Item(11, 2).Value = '[2013.0, 2.7267, 16.53669, 4.44]'
$ws.Cells.Item(11, 3).Value = '[2013.0, 2.7267, 16.53669, 4.44]'
$ws.Cells.Item(11, 4).Value = 'PASS'
$ws.Cells.Item(11, 5).Value = '[]'
$ws.Cells.Item(12, 1).Value = '(''AD37 7205 5893 V9N5 ODHY NA15'', ''Graphical User Interface'')'
$ws.Cells.Item(12, 2).Value = '[2007.0, 21.65938, 20.92826, 6.98]'
$ws.Cells.Item(12, 3).Value = '[2007.0, 21.65938, 20.92826, 6.98]'
$ws.Cells.Item(12, 4).Value = 'PASS'
$ws.Cells.Item(12, 5).Value = '[]'
$ws.Cells.Item(13, 1).Value = '(''AZ18 BRND WYIG EA0Y 06HS ZKNA TOEB'', ''service-desk'')'
$ws.Cells.Item(13, 2).Value = '[2004.0, 9.57692, 11.97564, 16.48]'
$ws.Cells.Item(13, 3).Value = '[2004.0, 9.57692, 11.97564, 16.48]'
$ws.Cells.Item(13, 4).Value = 'PASS'
$ws.Cells.Item(13, 5).Value = '[]'
$ws.Cells.Item(14, 1).Value = '(''MD66 NLWM 5VKM NHFK HNWO WCMH'', '''')'
$ws.Cells.Item(14, 2).Value = '[2018.0, 4.2544, 31.78594, 2.45]'
$ws.Cells.Item(14, 3).Value = '[2018.0, 4.2544, 31.78594, 2.45]'
$ws.Cells.Item(14, 4).Value = 'PASS'
$ws.Cells.Item(14, 5).Value = '[]'
$ws.Cells.Item(15, 1).Value = '(''AD49 5880 6175 MWHH EZC2 EPFK'', '''')'
$ws.Cells.Item(15, 2).Value = '[2016.0, 11.6998, 5.04259, 5.02]'
$ws.Cells.Item(15, 3).Value = '[2016.0, 11.6998, 5.04259, 5.02]'
$ws.Cells.Item(15, 4).Value = 'PASS'
$ws.Cells.Item(15, 5).Value = '[]'
$ws.Cells.Item(16, 1).Value = '(''CR10 1301 0843 5149 3134 7'', '''')'
$ws.Cells.Item(16, 2).Value = '[1996.0, 12.88227, 2.28652, 15.53]'
$ws.Cells.Item(16, 3).Value = '[1996.0, 12.88227, 2.28652, 15.53]'
$ws.Cells.Item(16, 4).Value = 'PASS'
$ws.Cells.Item(16, 5).Value = '[]'
$ws.Cells.Item(17, 1).Value = '(''CZ90 6388 1158 9085 1122 3730'', ''knowledge user'')'
$ws.Cells.Item(17, 2).Value = '[2011.0, 12.29505, 33.03468, 6.11]'
$ws.Cells.Item(17, 3).Value = '[2011.0, 12.29505, 33.03468, 6.11421]'
$ws.Cells.Item(17, 4).Value = 'FAIL'
$ws.Cells.Item(17, 5).Value = '[0.00421]'
$ws.Cells.Item(18, 1).Value = '(''LT17 0327 5423 3376 4902'', '''')'
$ws.Cells.Item(18, 2).Value = '[2013.0, 11.71432, 8.51316, 26.97]'
$ws.Cells.Item(18, 3).Value = '[2013.0, 11.71432, 8.51316, 26.97]'
$ws.Cells.Item(18, 4).Value = 'PASS'
$ws.Cells.Item(18, 5).Value = '[]'
$ws.Cells.Item(19, 1).Value = '(''MK21 288Y D1FL ESNP Q29'', '''', ''radical'')'
$ws.Cells.Item(19, 2).Value = '[2008.0, 23.93246, 33.43]'
$ws.Cells.Item(19, 3).Value = '[2008.0, 23.93246, 33.43]'
$ws.Cells.Item(19, 4).Value = 'PASS'
$ws.Cells.Item(19, 5).Value = '[]'
$ws.Cells.Item(20, 1).Value = '(''PL70 3380 5531 5257 6545 0735 9929'', ''Multi-layered'')'
$ws.Cells.Item(20, 2).Value = '[2019.0, 1.75088, 15.06291, 8.44]'
$ws.Cells.Item(20, 3).Value = '[2019.0, 1.75088, 15.06291, 78.44]'
$ws.Cells.Item(20, 4).Value = 'FAIL'
$ws.Cells.Item(20, 5).Value = '[70.0]'
$ws.Cells.Item(21, 1).Value = '(''HR34 1036 8633 1429 9393 2'', '''', '''')'
$ws.Cells.Item(21, 2).Value = '[2018.0, 32.0739, 15.4]'
$ws.Cells.Item(21, 3).Value = '[2018.0, 32.0739, 15.4]'
$ws.Cells.Item(21, 4).Value = 'PASS'
$ws.Cells.Item(21, 5).Value = '[]'
$ws.Cells.Item(22, 1).Value = '(''KZ45 2476 JXRW ABRQ 1XTC'', '''')'
$ws.Cells.Item(22, 2).Value = '[2014.0, 12.37228, 8.68935, 11.67]'
$ws.Cells.Item(22, 3).Value = '[2014.0, 12.37228, 8.68935, 11.67]'
$ws.Cells.Item(22, 4).Value = 'PASS'
$ws.Cells.Item(22, 5).Value = '[]'
$ws.Cells.Item(23, 1).Value = '(''AT32 9889 2449 9236 2133'', '''')'
$ws.Cells.Item(23, 2).Value = '[2014.0, 4.21371, 13.74512, 13.78]'
$ws.Cells.Item(23, 3).Value = '[2014.0, 4.21371, 13.74512, 13.78]'
$ws.Cells.Item(23, 4).Value = 'PASS'
$ws.Cells.Item(23, 5).Value = '[]'
$ws.Cells.Item(24, 1).Value = '(''DO91 L62F 0944 6935 1417 5974 5533'', ''solution-oriented'')'
$ws.Cells.Item(24, 2).Value = '[2011.0, 33.92766, 31.92266, 26.77]'
$ws.Cells.Item(24, 3).Value = '[2011.0, 33.92766, 31.92266, 26.77]'
$ws.Cells.Item(24, 4).Value = 'PASS'
$ws.Cells.Item(24, 5).Value = '[]'
$ws.Cells.Item(25, 1).Value = '(''AZ37 ZWXD 4WYS TYQI DDPB LZHH EMQS'', '''')'
$ws.Cells.Item(25, 2).Value = '[2021.0, 18.99396, 3.25298, 34.29]'
$ws.Cells.Item(25, 3).Value = '[2021.0, 18.99396, 3.25298, 34.29]'
$ws.Cells.Item(25, 4).Value = 'PASS'
$ws.Cells.Item(25, 5).Value = '[]'
$ws.Cells.Item(26, 1).Value = '(''FR71 1160 5711 20WC TLWU YZZY 355'', ''fresh-thinking'')'
$ws.Cells.Item(26, 2).Value = '[2009.0, 32.33383, 19.61646, 16.47]'
$ws.Cells.Item(26, 3).Value = '[2009.0, 32.33383, 19.61646, 16.47]'
$ws.Cells.Item(26, 4).Value = 'PASS'
$ws.Cells.Item(26, 5).Value = '[]'
$ws.Cells.Item(27, 1).Value = '(''SM93 B374 1481 278R NGOP UUZN R9P'', ''Compatible'')'
$ws.Cells.Item(27, 2).Value = '[2011.0, 15.29257, 24.09298, 11.7]'
$ws.Cells.Item(27, 3).Value = '[2011.0, 15.29257, 24.09298, 11.7]'
$ws.Cells.Item(27, 4).Value = 'PASS'
$ws.Cells.Item(27, 5).Value = '[]'
$ws.Cells.Item(28, 1).Value = '(''IS40 9180 4001 1378 4881 1907 39'', '''')'
$ws.Cells.Item(28, 2).Value = '[2009.0, 31.19619, 33.19297, 34.65]'
$ws.Cells.Item(28, 3).Value = '[2009.0, 31.19619, 33.19297, 34.65]'
$ws.Cells.Item(28, 4).Value = 'PASS'
$ws.Cells.Item(28, 5).Value = '[]'
$ws.Cells.Item(29, 1).Value = '(''GB84 DMNZ 2243 3445 8251 77'', ''Sharable'')'
$ws.Cells.Item(29, 2).Value = '[2004.0, 18.33324, 9.56963, 23.59]'
$ws.Cells.Item(29, 3).Value = '[2004.0, 18.33324, 9.56963, 23.59]'
$ws.Cells.Item(29, 4).Value = 'PASS'
$ws.Cells.Item(29, 5).Value = '[]'
$ws.Cells.Item(30, 1).Value = '(''LI45 5774 5DFL SQJK BW4N X'', ''projection'')'
$ws.Cells.Item(30, 2).Value = '[1999.0, 13.4668, 29.31703, 6.29]'
$ws.Cells.Item(30, 3).Value = '[1999.0, 13.4668, 29.31703, 6.29]'
$ws.Cells.Item(30, 4).Value = 'PASS'
$ws.Cells.Item(30, 5).Value = '[]'
$ws.Cells.Item(31, 1).Value = '(''NL51 SVWQ 7292 0904 45'', '''')'
$ws.Cells.Item(31, 2).Value = '[2015.0, 18.14148, 4.36244, 2.33]'
$ws.Cells.Item(31, 3).Value = '[2015.0, 18.14148, 4.36244, 2.33]'
$ws.Cells.Item(31, 4).Value = 'PASS'
$ws.Cells.Item(31, 5).Value = '[]'
$ws.Cells.Item(32, 1).Value = '(''MR71 3046 8079 7692 1026 2038 944'', '''')'
$ws.Cells.Item(32, 2).Value = '[2010.0, 11.52491, 3.78871, 21.97]'
$ws.Cells.Item(32, 3).Value = '[2010.0, 11.52491, 3.78871, 21.97]'
$ws.Cells.Item(32, 4).Value = 'PASS'
$ws.Cells.Item(32, 5).Value = '[]'
$ws.Cells.Item(33, 1).Value = '(''AL64 6952 6272 RTRO TG3H ZK4T 57UU'', ''needs-based'')'
$ws.Cells.Item(33, 2).Value = '[2017.0, 12.10915, 14.94297, 3.93]'
$ws.Cells.Item(33, 3).Value = '[2017.0, 12.10915, 14.94297, 3.93]'
$ws.Cells.Item(33, 4).Value = 'PASS'
$ws.Cells.Item(33, 5).Value = '[]'
$ws.Cells.Item(34, 1).Value = '(''BE79 3912 1524 3321'', '''')'
$ws.Cells.Item(34, 2).Value = '[0.0, 3.19658, 28.23225, 29.39]'
$ws.Cells.Item(34, 3).Value = '[0.0, 3.19658, 28.23225, 29.39]'
$ws.Cells.Item(34, 4).Value = 'PASS'
$ws.Cells.Item(34, 5).Value = '[]'
$ws.Cells.Item(35, 1).Value = '(''FI78 1545 2490 0436 58'', '''', '''')'
$ws.Cells.Item(35, 2).Value = '[2008.0, 19.87694, 19.17]'
$ws.Cells.Item(35, 3).Value = '[2008.0, 19.87694, 19.17]'
$ws.Cells.Item(35, 4).Value = 'PASS'
$ws.Cells.Item(35, 5).Value = '[]'
$ws.Cells.Item(36, 1).Value = '(''GR89 8461 8453 XZOY ULUP ZUHB 45H'', '''')'
$ws.Cells.Item(36, 2).Value = '[2016.0, 27.31544, 6.59009, 18.33]'
$ws.Cells.Item(36, 3).Value = '[2016.0, 27.31544, 6.59009, 18.33]'
$ws.Cells.Item(36, 4).Value = 'PASS'
$ws.Cells.Item(36, 5).Value = '[]'
$ws.Cells.Item(37, 1).Value = '(''FR44 5550 2233 12ZL XZL2 6UMG N98'', ''Face to face'')'
$ws.Cells.Item(37, 2).Value = '[2007.0, 14.8618, 19.56054, 2.98]'
$ws.Cells.Item(37, 3).Value = '[2007.0, 14.8618, 19.56054, 2.98]'
$ws.Cells.Item(37, 4).Value = 'PASS'
$ws.Cells.Item(37, 5).Value = '[]'
$ws.Cells.Item(38, 1).Value = '(''ME97 4838 3055 2140 4652 40'', '''')'
$ws.Cells.Item(38, 2).Value = '[2009.0, 15.54734, 5.07709, 34.93]'
$ws.Cells.Item(38, 3).Value = '[2009.0, 15.54734, 5.07709, 34.93]'
$ws.Cells.Item(38, 4).Value = 'PASS'
$ws.Cells.Item(38, 5).Value = '[]'
$ws.Cells.Item(39, 1).Value = '(''FR45 0872 9704 65YG MLW0 Y1RJ I95'', '''')'
$ws.Cells.Item(39, 2).Value = '[2014.0, 25.56047, 2.36789, 32.25]'
$ws.Cells.Item(39, 3).Value = '[2014.0, 25.56047, 2.36789, 32.25]'
$ws.Cells.Item(39, 4).Value = 'PASS'
$ws.Cells.Item(39, 5).Value = '[]'
$ws.Cells.Item(40, 1).Value = '(''RO10 PHVB YLH0 30YU EB1K NQEE'', '''')'
$ws.Cells.Item(40, 2).Value = '[2014.0, 1.7785, 4.9239, 33.51]'
$ws.Cells.Item(40, 3).Value = '[2014.0, 1.7785, 4.9239, 33.51]'
$ws.Cells.Item(40, 4).Value = 'PASS'
$ws.Cells.Item(40, 5).Value = '[]'
$ws.Cells.Item(41, 1).Value = '(''PT16 5670 0074 0937 5625 5297 5'', ''help-desk'')'
$ws.Cells.Item(41, 2).Value = '[2018.0, 29.02928, 19.03125, 29.01]'
$ws.Cells.Item(41, 3).Value = '[2018.0, 29.02928, 19.03125, 29.01]'
$ws.Cells.Item(41, 4).Value = 'PASS'
$ws.Cells.Item(41, 5).Value = '[]'
$ws.Cells.Item(42, 1).Value = '(''SM53 L335 2136 439R IEUF WEC5 ZRD'', ''data-warehouse'')'
$ws.Cells.Item(42, 2).Value = '[2014.0, 26.48204, 18.89291, 21.14]'
$ws.Cells.Item(42, 3).Value = '[2014.0, 26.48204, 18.89291, 21.14]'
$ws.Cells.Item(42, 4).Value = 'PASS'
$ws.Cells.Item(42, 5).Value = '[]'
$ws.Cells.Item(43, 1).Value = '(''LB70 9084 U7AE PV68 JRC9 SLPY OAZA'', ''Right-sized'')'
$ws.Cells.Item(43, 2).Value = '[2015.0, 31.51073, 28.31371, 3.38]'
$ws.Cells.Item(43, 3).Value = '[2015.0, 31.51073, 28.31371, 3.38]'
$ws.Cells.Item(43, 4).Value = 'PASS'
$ws.Cells.Item(43, 5).Value = '[]'
$ws.Cells.Item(44, 1).Value = '(''BE13 3895 5129 7876'', ''Switchable'')'
$ws.Cells.Item(44, 2).Value = '[2004.0, 6.39977, 30.0622, 10.57]'
$ws.Cells.Item(44, 3).Value = '[2004.0, 6.39977, 30.0622, 10.57]'
$ws.Cells.Item(44, 4).Value = 'PASS'
$ws.Cells.Item(44, 5).Value = '[]'
$ws.Cells.Item(45, 1).Value = '(''SA98 85XA 647W JRQ2 V4CK URMJ'', ''implementation'')'
$ws.Cells.Item(45, 2).Value = '[2015.0, 4.65345, 13.77682, 13.27]'
$ws.Cells.Item(45, 3).Value = '[2015.0, 4.65345, 13.77682, 13.27]'
$ws.Cells.Item(45, 4).Value = 'PASS'
$ws.Cells.Item(45, 5).Value = '[]'
$ws.Cells.Item(46, 1).Value = '(''PS08 LFNK Z2L2 IZBZ TATP PMIE 0ORO T'', '''')'
$ws.Cells.Item(46, 2).Value = '[2013.0, 31.74558, 21.62143, 10.93]'
$ws.Cells.Item(46, 3).Value = '[2013.0, 31.74558, 21.62143, 10.93]'
$ws.Cells.Item(46, 4).Value = 'PASS'
$ws.Cells.Item(46, 5).Value = '[]'
$ws.Cells.Item(47, 1).Value = '(''HR23 9970 5212 1133 4458 3'', ''array'')'
$ws.Cells.Item(47, 2).Value = '[2009.0, 8.59499, 31.11391, 18.98]'
$ws.Cells.Item(47, 3).Value = '[2009.0, 8.59499, 31.11391, 18.98]'
$ws.Cells.Item(47, 4).Value = 'PASS'
$ws.Cells.Item(47, 5).Value = '[]'
$ws.Cells.Item(48, 1).Value = '(''LB25 8193 ZZVU ZNFI UKTG CIWA JSQV'', '''')'
$ws.Cells.Item(48, 2).Value = '[2011.0, 33.64968, 25.22065, 9.9]'
$ws.Cells.Item(48, 3).Value = '[2011.0, 33.64968, 25.22065, 9.9]'
$ws.Cells.Item(48, 4).Value = 'PASS'
$ws.Cells.Item(48, 5).Value = '[]'
$ws.Cells.Item(49, 1).Value = '(''FR59 2632 5348 89TN HTWN GYRY C99'', '''')'
$ws.Cells.Item(49, 2).Value = '[2011.0, 30.4217, 11.26172, 14.88]'
$ws.Cells.Item(49, 3).Value = '[2011.0, 30.4217, 11.26172, 14.88]'
$ws.Cells.Item(49, 4).Value = 'PASS'
$ws.Cells.Item(49, 5).Value = '[]'
$ws.Cells.Item(50, 1).Value = '(''TR21 9920 7ENO DPGN 3QKY 7XWB XG'', '''')'
$ws.Cells.Item(50, 2).Value = '[2009.0, 14.85743, 24.29145, 17.76]'
$ws.Cells.Item(50, 3).Value = '[2009.0, 14.85743, 24.29145, 17.76]'
$ws.Cells.Item(50, 4).Value = 'PASS'
$ws.Cells.Item(50, 5).Value = '[]'
$ws.Cells.Item(51, 1).Value = '(''IS51 2444 7567 4951 4366 0772 89'', '''')'
$ws.Cells.Item(51, 2).Value = '[2005.0, 19.12235, 34.37349, 13.68]'
$ws.Cells.Item(51, 3).Value = '[2005.0, 19.12235, 34.37349, 13.68]'
$ws.Cells.Item(51, 4).Value = 'PASS'
$ws.Cells.Item(51, 5).Value = '[]'
$ws.Cells.Item(52, 1).Value = '(''SI47 6748 4442 3463 823'', '''')'
$ws.Cells.Item(52, 2).Value = '[1998.0, 18.09311, 5.76609, 2.06]'
$ws.Cells.Item(52, 3).Value = '[1998.0, 18.09311, 5.76609, 2.06]'
$ws.Cells.Item(52, 4).Value = 'PASS'
$ws.Cells.Item(52, 5).Value = '[]'
$ws.Cells.Item(53, 1).Value = '(''SI68 1820 7430 9240 077'', '''')'
$ws.Cells.Item(53, 2).Value = '[2018.0, 7.37854, 19.69095, 20.21]'
$ws.Cells.Item(53, 3).Value = '[2018.0, 7.37854, 19.69095, 20.21]'
$ws.Cells.Item(53, 4).Value = 'PASS'
$ws.Cells.Item(53, 5).Value = '[]'
$ws.Cells.Item(54, 1).Value = '(''PS25 EPKD CHVW OANF DW5Q RN0W WTCX K'', '''')'
$ws.Cells.Item(54, 2).Value = '[2019.0, 8.35226, 3.71473, 11.16]'
$ws.Cells.Item(54, 3).Value = '[2019.0, 8.35226, 3.71473, 11.16]'
$ws.Cells.Item(54, 4).Value = 'PASS'
$ws.Cells.Item(54, 5).Value = '[]'
$ws.Cells.Item(55, 1).Value = '(''GE04 VU99 6596 7310 5692 42'', '''')'
$ws.Cells.Item(55, 2).Value = '[2018.0, 9.60983, 16.88953, 1.17]'
$ws.Cells.Item(55, 3).Value = '[2018.0, 9.60983, 16.88953, 1.17]'
$ws.Cells.Item(55, 4).Value = 'PASS'
$ws.Cells.Item(55, 5).Value = '[]'
$ws.Cells.Item(56, 1).Value = '(''AT26 5820 2675 0356 3271'', ''zero defect'')'
$ws.Cells.Item(56, 2).Value = '[2008.0, 28.49054, 14.7059, 27.72]'
$ws.Cells.Item(56, 3).Value = '[2008.0, 28.49054, 14.7059, 27.72]'
$ws.Cells.Item(56, 4).Value = 'PASS'
$ws.Cells.Item(56, 5).Value = '[]'
$ws.Cells.Item(57, 1).Value = '(''BG90 CKND 3854 18HU XD5M 8J'', '''')'
$ws.Cells.Item(57, 2).Value = '[2001.0, 4.03151, 2.32338, 12.23]'
$ws.Cells.Item(57, 3).Value = '[2001.0, 4.03151, 2.32338, 12.23]'
$ws.Cells.Item(57, 4).Value = 'PASS'
$ws.Cells.Item(57, 5).Value = '[]'
$ws.Cells.Item(58, 1).Value = '(''BG95 VNYO 7604 271R NNAK MA'', '''')'
$ws.Cells.Item(58, 2).Value = '[2016.0, 13.51647, 13.46117, 3.87]'
$ws.Cells.Item(58, 3).Value = '[2016.0, 13.51647, 13.46117, 3.87]'
$ws.Cells.Item(58, 4).Value = 'PASS'
$ws.Cells.Item(58, 5).Value = '[]'
$ws.Cells.Item(59, 1).Value = '(''HR32 2437 0990 3468 3251 7'', '''', '''')'
$ws.Cells.Item(59, 2).Value = '[2016.0, 7.17604, 21.05]'
$ws.Cells.Item(59, 3).Value = '[2016.0, 7.17604, 21.05]'
$ws.Cells.Item(59, 4).Value = 'PASS'
$ws.Cells.Item(59, 5).Value = '[]'
$ws.Cells.Item(60, 1).Value = '(''MU73 TSFB 6869 9727 7482 6943 284V LT'', '''')'
$ws.Cells.Item(60, 2).Value = '[2012.0, 16.71725, 18.01337, 9.79]'
$ws.Cells.Item(60, 3).Value = '[2012.0, 16.71725, 18.01337, 9.79]'
$ws.Cells.Item(60, 4).Value = 'PASS'
$ws.Cells.Item(60, 5).Value = '[]'
$ws.Cells.Item(61, 1).Value = '(''MC48 7184 8197 768V CY1K IMGJ R14'', '''', '''')'
$ws.Cells.Item(61, 2).Value = '[2016.0, 21.33205, 16.03]'
$ws.Cells.Item(61, 3).Value = '[2016.0, 21.33205, 16.03]'
$ws.Cells.Item(61, 4).Value = 'PASS'
$ws.Cells.Item(61, 5).Value = '[]'
$ws.Cells.Item(62, 1).Value = '(''GR16 3065 793B Z5RU G17Q Q3QN RXW'', '''')'
$ws.Cells.Item(62, 2).Value = '[2009.0, 3.85894, 2.76212, 26.07]'
$ws.Cells.Item(62, 3).Value = '[2009.0, 3.85894, 2.76212, 26.07]'
$ws.Cells.Item(62, 4).Value = 'PASS'
$ws.Cells.Item(62, 5).Value = '[]'
$ws.Cells.Item(63, 1).Value = '(''GL19 8754 3034 5493 70'', ''Cloned'')'
$ws.Cells.Item(63, 2).Value = '[2012.0, 27.3021, 19.08725, 19.51]'
$ws.Cells.Item(63, 3).Value = '[2012.0, 27.3021, 19.08725, 19.51]'
$ws.Cells.Item(63, 4).Value = 'PASS'
$ws.Cells.Item(63, 5).Value = '[]'
$ws.Cells.Item(64, 1).Value = '(''MR63 4626 5538 5305 3799 0332 886'', ''mission-critical'')'
$ws.Cells.Item(64, 2).Value = '[2015.0, 11.05505, 23.72029, 5.36]'
$ws.Cells.Item(64, 3).Value = '[2015.0, 11.05505, 23.72029, 5.36]'
$ws.Cells.Item(64, 4).Value = 'PASS'
$ws.Cells.Item(64, 5).Value = '[]'
$ws.Cells.Item(65, 1).Value = '(''LV97 HJQI XSFX MBNZ TMMA N'', ''moratorium'')'
$ws.Cells.Item(65, 2).Value = '[2018.0, 8.78126, 3.72063, 2.56]'
$ws.Cells.Item(65, 3).Value = '[2018.0, 8.78126, 3.72063, 2.56]'
$ws.Cells.Item(65, 4).Value = 'PASS'
$ws.Cells.Item(65, 5).Value = '[]'
$ws.Cells.Item(66, 1).Value = '(''DK13 2114 5766 0520 95'', '''')'
$ws.Cells.Item(66, 2).Value = '[2019.0, 20.32835, 22.20597, 22.99]'
$ws.Cells.Item(66, 3).Value = '[2019.0, 20.32835, 22.20597, 22.99]'
$ws.Cells.Item(66, 4).Value = 'PASS'
$ws.Cells.Item(66, 5).Value = '[]'
$ws.Cells.Item(67, 1).Value = '(''AD89 4822 4731 3CVH WNRZ ONKG'', '''')'
$ws.Cells.Item(67, 2).Value = '[2020.0, 13.15791, 27.04429, 30.74]'
$ws.Cells.Item(67, 3).Value = '[2020.0, 13.15791, 27.04429, 30.74]'
$ws.Cells.Item(67, 4).Value = 'PASS'
$ws.Cells.Item(67, 5).Value = '[]'
$ws.Cells.Item(68, 1).Value = '(''SI74 2819 6967 0008 990'', ''Open-architected'')'
$ws.Cells.Item(68, 2).Value = '[2013.0, 31.38258, 19.83445, 3.91]'
$ws.Cells.Item(68, 3).Value = '[2013.0, 31.38258, 19.83445, 3.91]'
$ws.Cells.Item(68, 4).Value = 'PASS'
$ws.Cells.Item(68, 5).Value = '[]'
$ws.Cells.Item(69, 1).Value = '(''VG31 ZQOE 1664 3145 5498 7770'', ''Decentralized'')'
$ws.Cells.Item(69, 2).Value = '[2017.0, 26.99993, 19.80792, 33.3]'
$ws.Cells.Item(69, 3).Value = '[2017.0, 26.99993, 19.80792, 33.3]'
$ws.Cells.Item(69, 4).Value = 'PASS'
$ws.Cells.Item(69, 5).Value = '[]'
$ws.Cells.Item(70, 1).Value = '(''HU31 3766 8395 6898 7567 1070 2083'', '''')'
$ws.Cells.Item(70, 2).Value = '[2001.0, 18.06291, 32.67256, 1.26]'
$ws.Cells.Item(70, 3).Value = '[2001.0, 18.06291, 32.67256, 1.26]'
$ws.Cells.Item(70, 4).Value = 'PASS'
$ws.Cells.Item(70, 5).Value = '[]'
$ws.Cells.Item(71, 1).Value = '(''LU02 283I ZHUB SALT NHTF'', '''')'
$ws.Cells.Item(71, 2).Value = '[2012.0, 6.07496, 16.91051, 8.29]'
$ws.Cells.Item(71, 3).Value = '[2012.0, 6.07496, 16.91051, 8.29]'
$ws.Cells.Item(71, 4).Value = 'PASS'
$ws.Cells.Item(71, 5).Value = '[]'
$ws.Cells.Item(72, 1).Value = '(''LB69 6641 Z6AV VCAK 7SZB FHJY TXDY'', '''')'
$ws.Cells.Item(72, 2).Value = '[2005.0, 31.62293, 23.62821, 20.67]'
$ws.Cells.Item(72, 3).Value = '[2005.0, 31.62293, 23.62821, 20.67]'
$ws.Cells.Item(72, 4).Value = 'PASS'
$ws.Cells.Item(72, 5).Value = '[]'
$ws.Cells.Item(73, 1).Value = '(''SE62 8794 0081 7723 6087 8593'', '''')'
$ws.Cells.Item(73, 2).Value = '[2016.0, 23.16398, 30.78324, 32.55]'
$ws.Cells.Item(73, 3).Value = '[2016.0, 23.16398, 30.78324, 32.55]'
$ws.Cells.Item(73, 4).Value = 'PASS'
$ws.Cells.Item(73, 5).Value = '[]'
$ws.Cells.Item(74, 1).Value = '(''FR55 8614 4500 58LJ 9ONQ WO5M V56'', ''Object-based'')'
$ws.Cells.Item(74, 2).Value = '[2010.0, 8.77452, 29.77394, 28.4]'
$ws.Cells.Item(74, 3).Value = '[2010.0, 8.77452, 29.77394, 28.4]'
$ws.Cells.Item(74, 4).Value = 'PASS'
$ws.Cells.Item(74, 5).Value = '[]'
$ws.Cells.Item(75, 1).Value = '(''IS96 7839 4200 8268 2110 4652 56'', ''firmware'')'
$ws.Cells.Item(75, 2).Value = '[2007.0, 2.86685, 22.97517, 19.45]'
$ws.Cells.Item(75, 3).Value = '[2052.0, 2.86685, 22.97517, 19.45]'
$ws.Cells.Item(75, 4).Value = 'FAIL'
$ws.Cells.Item(75, 5).Value = '[45.0]'
$ws.Cells.Item(76, 1).Value = '(''FR75 2140 3789 37VT NT7A FUCF P92'', ''explicit'')'
$ws.Cells.Item(76, 2).Value = '[2020.0, 1.16374, 33.73684, 22.43]'
$ws.Cells.Item(76, 3).Value = '[2020.0, 1.16374, 33.73684, 22.43]'
$ws.Cells.Item(76, 4).Value = 'PASS'
$ws.Cells.Item(76, 5).Value = '[]'
$ws.Cells.Item(77, 1).Value = '(''PS08 QCPE USPO K7TA E3FQ GQOQ KRF3 T'', '''')'
$ws.Cells.Item(77, 2).Value = '[2018.0, 15.52732, 14.93181, 21.89]'
$ws.Cells.Item(77, 3).Value = '[2018.0, 15.52732, 14.93181, 21.89]'
$ws.Cells.Item(77, 4).Value = 'PASS'
$ws.Cells.Item(77, 5).Value = '[]'
$ws.Cells.Item(78, 1).Value = '(''MK46 734T CVK7 JKIO B11'', '''')'
$ws.Cells.Item(78, 2).Value = '[0.0, 17.17025, 10.10083, 7.62]'
$ws.Cells.Item(78, 3).Value = '[0.0, 17.17025, 10.10083, 7.62]'
$ws.Cells.Item(78, 4).Value = 'PASS'
$ws.Cells.Item(78, 5).Value = '[]'
$ws.Cells.Item(79, 1).Value = '(''FR45 8691 9784 71NY HOEF 55CI N94'', '''')'
$ws.Cells.Item(79, 2).Value = '[2012.0, 4.53097, 31.62524, 19.08]'
$ws.Cells.Item(79, 3).Value = '[2012.0, 4.53097, 31.62524, 19.08]'
$ws.Cells.Item(79, 4).Value = 'PASS'
$ws.Cells.Item(79, 5).Value = '[]'
$ws.Cells.Item(80, 1).Value = '(''MT93 PUXX 1605 7YXT YCML KSPR 0UNT 3IH'', '''')'
$ws.Cells.Item(80, 2).Value = '[2014.0, 15.13911, 5.78656, 19.52]'
$ws.Cells.Item(80, 3).Value = '[2014.0, 15.13911, 5.78656, 19.52]'
$ws.Cells.Item(80, 4).Value = 'PASS'
$ws.Cells.Item(80, 5).Value = '[]'
$ws.Cells.Item(81, 1).Value = '(''AE58 6169 4300 9504 9513 851'', '''')'
$ws.Cells.Item(81, 2).Value = '[2005.0, 16.79262, 34.22827, 8.82]'
$ws.Cells.Item(81, 3).Value = '[2005.0, 16.79262, 34.22827, 8.82]'
$ws.Cells.Item(81, 4).Value = 'PASS'
$ws.Cells.Item(81, 5).Value = '[]'
$ws.Cells.Item(82, 1).Value = '(''FR75 1304 6692 36HH YNNW EDRP U16'', '''')'
$ws.Cells.Item(82, 2).Value = '[2017.0, 11.65434, 15.00641, 27.21]'
$ws.Cells.Item(82, 3).Value = '[2017.0, 11.65434, 15.00641, 27.21]'
$ws.Cells.Item(82, 4).Value = 'PASS'
$ws.Cells.Item(82, 5).Value = '[]'
$ws.Cells.Item(83, 1).Value = '(''SI70 8014 5983 9508 969'', '''')'
$ws.Cells.Item(83, 2).Value = '[2011.0, 20.02967, 16.85391, 21.54]'
$ws.Cells.Item(83, 3).Value = '[2011.0, 20.02967, 16.85391, 21.54]'
$ws.Cells.Item(83, 4).Value = 'PASS'
$ws.Cells.Item(83, 5).Value = '[]'
$ws.Cells.Item(84, 1).Value = '(''CZ91 9143 2822 8759 8003 0329'', '''')'
$ws.Cells.Item(84, 2).Value = '[2020.0, 30.56008, 14.6595, 9.0]'
$ws.Cells.Item(84, 3).Value = '[2020.0, 30.56008, 14.6595, 9.0]'
$ws.Cells.Item(84, 4).Value = 'PASS'
$ws.Cells.Item(84, 5).Value = '[]'
$ws.Cells.Item(85, 1).Value = '(''MD23 PZGO HHRM ARCT VVLJ JRGC'', ''project'')'
$ws.Cells.Item(85, 2).Value = '[0.0, 24.61024, 21.53354, 25.02]'
$ws.Cells.Item(85, 3).Value = '[0.0, 24.61024, 21.53354, 25.02]'
$ws.Cells.Item(85, 4).Value = 'PASS'
$ws.Cells.Item(85, 5).Value = '[]'
$ws.Cells.Item(86, 1).Value = '(''FR06 5787 2577 99HY OXFG H5AC Y13'', ''conglomeration'')'
$ws.Cells.Item(86, 2).Value = '[2014.0, 17.54223, 3.39843, 11.09]'
$ws.Cells.Item(86, 3).Value = '[2014.0, 17.54223, 3.39843, 11.09]'
$ws.Cells.Item(86, 4).Value = 'PASS'
$ws.Cells.Item(86, 5).Value = '[]'
$ws.Cells.Item(87, 1).Value = '(''FR14 9721 1762 58D9 L9BL TAW2 Z76'', ''encoding'')'
$ws.Cells.Item(87, 2).Value = '[2007.0, 31.63965, 10.49419, 20.16]'
$ws.Cells.Item(87, 3).Value = '[2007.0, 31.63965, 10.49419, 20.16]'
$ws.Cells.Item(87, 4).Value = 'PASS'
$ws.Cells.Item(87, 5).Value = '[]'
$ws.Cells.Item(88, 1).Value = '(''BA56 6184 0204 5174 8421'', '''')'
$ws.Cells.Item(88, 2).Value = '[2015.0, 14.56754, 34.13968, 21.32]'
$ws.Cells.Item(88, 3).Value = '[2015.0, 14.56754, 34.13968, 21.32]'
$ws.Cells.Item(88, 4).Value = 'PASS'
$ws.Cells.Item(88, 5).Value = '[]'
$ws.Cells.Item(89, 1).Value = '(''BH52 GMXD DQGK BQV8 QYZO QW'', '''')'
$ws.Cells.Item(89, 2).Value = '[2010.0, 19.81187, 28.02982, 14.0]'
$ws.Cells.Item(89, 3).Value = '[2010.0, 19.81187, 28.02982, 14.0]'
$ws.Cells.Item(89, 4).Value = 'PASS'
$ws.Cells.Item(89, 5).Value = '[]'
$ws.Cells.Item(90, 1).Value = '(''LV46 XXIX NMFQ 0BWN U3ZD M'', '''')'
$ws.Cells.Item(90, 2).Value = '[2009.0, 14.8385, 9.91532, 20.1]'
$ws.Cells.Item(90, 3).Value = '[2009.0, 14.8385, 9.91532, 20.1]'
$ws.Cells.Item(90, 4).Value = 'PASS'
$ws.Cells.Item(90, 5).Value = '[]'
$ws.Cells.Item(91, 1).Value = '(''CY19 2831 7159 XHH9 SPNB U3ZX H0G1'', '''')'
$ws.Cells.Item(91, 2).Value = '[2016.0, 20.68625, 13.45405, 18.54]'
$ws.Cells.Item(91, 3).Value = '[2016.0, 20.68625, 13.45405, 18.54]'
$ws.Cells.Item(91, 4).Value = 'PASS'
$ws.Cells.Item(91, 5).Value = '[]'
$ws.Cells.Item(92, 1).Value = '(''SA47 78DR HSF9 1PO2 VXQC KBL3'', '''')'
$ws.Cells.Item(92, 2).Value = '[2016.0, 8.04129, 18.91455, 30.66]'
$ws.Cells.Item(92, 3).Value = '[2016.0, 8.04129, 18.91455, 30.66]'
$ws.Cells.Item(92, 4).Value = 'PASS'
$ws.Cells.Item(92, 5).Value = '[]'
$ws.Cells.Item(93, 1).Value = '(''NL03 DJJO 7107 5717 84'', '''')'
$ws.Cells.Item(93, 2).Value = '[2016.0, 26.35799, 2.4478, 21.3]'
$ws.Cells.Item(93, 3).Value = '[2016.0, 26.35799, 2.4478, 21.3]'
$ws.Cells.Item(93, 4).Value = 'PASS'
$ws.Cells.Item(93, 5).Value = '[]'
$ws.Cells.Item(94, 1).Value = '(''PK46 ZUUK KNR1 DKPJ 31JF A2AU'', '''')'
$ws.Cells.Item(94, 2).Value = '[2014.0, 26.11139, 19.87706, 19.74]'
$ws.Cells.Item(94, 3).Value = '[2014.0, 26.11139, 19.87706, 19.74]'
$ws.Cells.Item(94, 4).Value = 'PASS'
$ws.Cells.Item(94, 5).Value = '[]'
$ws.Cells.Item(95, 1).Value = '(''FR17 4386 6057 70VE CKFI L4ST U61'', '''')'
$ws.Cells.Item(95, 2).Value = '[2016.0, 10.72207, 26.59578, 27.96]'
$ws.Cells.Item(95, 3).Value = '[2016.0, 10.72207, 26.59578, 27.96]'
$ws.Cells.Item(95, 4).Value = 'PASS'
$ws.Cells.Item(95, 5).Value = '[]'
$ws.Cells.Item(96, 1).Value = '(''FR26 6978 0875 43UP MGMY PJRT 242'', ''content-based'')'
$ws.Cells.Item(96, 2).Value = '[2018.0, 29.333, 33.58855, 32.64]'
$ws.Cells.Item(96, 3).Value = '[2018.0, 29.333, 33.58855, 32.64]'
$ws.Cells.Item(96, 4).Value = 'PASS'
$ws.Cells.Item(96, 5).Value = '[]'
$ws.Cells.Item(97, 1).Value = '(''FR85 0875 9198 14NK HWJH HMHV 101'', '''')'
$ws.Cells.Item(97, 2).Value = '[2015.0, 3.16178, 15.56012, 20.09]'
$ws.Cells.Item(97, 3).Value = '[2015.0, 3.16178, 15.56012, 20.09]'
$ws.Cells.Item(97, 4).Value = 'PASS'
$ws.Cells.Item(97, 5).Value = '[]'
$ws.Cells.Item(98, 1).Value = '(''CR36 5946 7540 0249 7042 5'', ''Customer-focused'')'
$ws.Cells.Item(98, 2).Value = '[2010.0, 24.41216, 9.40856, 2.4]'
$ws.Cells.Item(98, 3).Value = '[2010.0, 24.41216, 9.40856, 2.4]'
$ws.Cells.Item(98, 4).Value = 'PASS'
$ws.Cells.Item(98, 5).Value = '[]'
$ws.Cells.Item(99, 1).Value = '(''MT62 RGJF 9551 4JYP 1U32 WSUA G82L BQJ'', '''')'
$ws.Cells.Item(99, 2).Value = '[2015.0, 30.01065, 32.08387, 24.83]'
$ws.Cells.Item(99, 3).Value = '[2015.0, 30.01065, 32.08387, 24.83]'
$ws.Cells.Item(99, 4).Value = 'PASS'
$ws.Cells.Item(99, 5).Value = '[]'
$ws.Cells.Item(100, 1).Value = '(''PS65 MCEE WBAV 8APE QGL0 Z1N4 QF8J 8'', '''')'
$ws.Cells.Item(100, 2).Value = '[2015.0, 27.52542, 31.74467, 22.4]'
$ws.Cells.Item(100, 3).Value = '[2015.0, 27.52542, 31.74467, 22.4]'
$ws.Cells.Item(100, 4).Value = 'PASS'
$ws.Cells.Item(100, 5).Value = '[]'
$ws.Cells.Item(101, 1).Value = '(''FR62 9845 5528 86SZ 00AZ EKBP C94'', '''')'
$ws.Cells.Item(101, 2).Value = '[2011.0, 4.02817, 10.08147, 2.2]'
$ws.Cells.Item(101, 3).Value = '[2011.0, 4.02817, 10.08147, 2.2]'
$ws.Cells.Item(101, 4).Value = 'PASS'
$ws.Cells.Item(101, 5).Value = '[]'
$ws.Cells.Item(102, 1).Value = '(''NL25 OCGO 6836 2647 49'', ''transitional'')'
$ws.Cells.Item(102, 2).Value = '[2012.0, 11.94561, 28.36059, 9.13]'
$ws.Cells.Item(102, 3).Value = '[2012.0, 11.94561, 28.36059, 9.13]'
$ws.Cells.Item(102, 4).Value = 'PASS'
$ws.Cells.Item(102, 5).Value = '[]'
$ws.Cells.Item(103, 1).Value = '(''PS73 PCWV 2NHC WKIA 8UXL URY4 AU7Y 3'', '''')'
$ws.Cells.Item(103, 2).Value = '[2012.0, 32.86778, 32.34208, 15.86]'
$ws.Cells.Item(103, 3).Value = '[2012.0, 32.86778, 32.34208, 15.86]'
$ws.Cells.Item(103, 4).Value = 'PASS'
$ws.Cells.Item(103, 5).Value = '[]'
$ws.Cells.Item(104, 1).Value = '(''LI09 7108 2V6S B11T BR1F R'', '''')'
$ws.Cells.Item(104, 2).Value = '[2015.0, 23.51244, 20.43548, 9.47]'
$ws.Cells.Item(104, 3).Value = '[2015.0, 23.51244, 20.43548, 9.47]'
$ws.Cells.Item(104, 4).Value = 'PASS'
$ws.Cells.Item(104, 5).Value = '[]'
$ws.Cells.Item(105, 1).Value = '(''AD25 5320 1695 V1J5 XQ74 I0R4'', ''standardization'')'
$ws.Cells.Item(105, 2).Value = '[2011.0, 13.12997, 5.41229, 2.24]'
$ws.Cells.Item(105, 3).Value = '[2011.0, 13.12997, 5.41229, 2.24]'
$ws.Cells.Item(105, 4).Value = 'PASS'
$ws.Cells.Item(105, 5).Value = '[]'
$ws.Cells.Item(106, 1).Value = '(''AT21 0815 5067 4654 1583'', '''')'
$ws.Cells.Item(106, 2).Value = '[2012.0, 25.21563, 5.92503, 11.05]'
$ws.Cells.Item(106, 3).Value = '[2012.0, 25.21563, 5.92503, 11.05]'
$ws.Cells.Item(106, 4).Value = 'PASS'
$ws.Cells.Item(106, 5).Value = '[]'
$ws.Cells.Item(107, 1).Value = '(''IL37 0716 3099 8553 0951 409'', '''')'
$ws.Cells.Item(107, 2).Value = '[2010.0, 24.14219, 13.83815, 17.75]'
$ws.Cells.Item(107, 3).Value = '[2010.0, 24.14219, 13.83815, 17.75]'
$ws.Cells.Item(107, 4).Value = 'PASS'
$ws.Cells.Item(107, 5).Value = '[]'
$ws.Cells.Item(108, 1).Value = '(''FR43 0281 5120 22UO EEV5 NNHE 937'', ''zero administration'')'
$ws.Cells.Item(108, 2).Value = '[2015.0, 1.44457, 23.31404, 9.38]'
$ws.Cells.Item(108, 3).Value = '[2015.0, 1.44457, 23.31404, 9.38]'
$ws.Cells.Item(108, 4).Value = 'PASS'
$ws.Cells.Item(108, 5).Value = '[]'
$ws.Cells.Item(109, 1).Value = '(''VG12 VRDR 7201 1116 1917 7386'', '''')'
$ws.Cells.Item(109, 2).Value = '[2020.0, 26.35827, 21.27747, 9.2]'
$ws.Cells.Item(109, 3).Value = '[2020.0, 26.35827, 21.27747, 9.2]'
$ws.Cells.Item(109, 4).Value = 'PASS'
$ws.Cells.Item(109, 5).Value = '[]'
$ws.Cells.Item(110, 1).Value = '(''AL72 8333 0306 154N G73I 7JSM 6BI1'', ''Integrated'')'
$ws.Cells.Item(110, 2).Value = '[2005.0, 16.76123, 23.15892, 21.91]'
$ws.Cells.Item(110, 3).Value = '[2005.0, 16.76123, 23.15892, 21.91]'
$ws.Cells.Item(110, 4).Value = 'PASS'
$ws.Cells.Item(110, 5).Value = '[]'
$ws.Cells.Item(111, 1).Value = '(''FR34 8487 2432 90RA D8D0 QFVP 865'', ''local area network'')'
$ws.Cells.Item(111, 2).Value = '[2018.0, 21.6696, 25.44197, 19.79]'
$ws.Cells.Item(111, 3).Value = '[2018.0, 21.6696, 25.44197, 19.79]'
$ws.Cells.Item(111, 4).Value = 'PASS'
$ws.Cells.Item(111, 5).Value = '[]'
$ws.Cells.Item(112, 1).Value = '(''BE85 6513 7708 8282'', '''')'
$ws.Cells.Item(112, 2).Value = '[2018.0, 20.28775, 34.65946, 13.52]'
$ws.Cells.Item(112, 3).Value = '[2018.0, 20.28775, 34.65946, 13.52]'
$ws.Cells.Item(112, 4).Value = 'PASS'
$ws.Cells.Item(112, 5).Value = '[]'
$ws.Cells.Item(113, 1).Value = '(''IT74 D323 7563 152T 29DG ADUC J46'', ''well-modulated'')'
$ws.Cells.Item(113, 2).Value = '[2015.0, 12.93085, 1.32315, 31.79]'
$ws.Cells.Item(113, 3).Value = '[2015.0, 12.93085, 1.32315, 31.79]'
$ws.Cells.Item(113, 4).Value = 'PASS'
$ws.Cells.Item(113, 5).Value = '[]'
$ws.Cells.Item(114, 1).Value = '(''LB16 5806 PU1J JVSB F9DV JKDI PXRM'', '''')'
$ws.Cells.Item(114, 2).Value = '[2004.0, 15.35258, 3.82453, 4.05]'
$ws.Cells.Item(114, 3).Value = '[2004.0, 15.35258, 3.82453, 4.05]'
$ws.Cells.Item(114, 4).Value = 'PASS'
$ws.Cells.Item(114, 5).Value = '[]'
$ws.Cells.Item(115, 1).Value = '(''AZ08 LFBO EUZA W4OP 9B1Y DT3H NSIN'', ''Object-based'')'
$ws.Cells.Item(115, 2).Value = '[2010.0, 21.78551, 3.16202, 6.33]'
$ws.Cells.Item(115, 3).Value = '[2010.0, 21.78551, 3.16202, 6.33]'
$ws.Cells.Item(115, 4).Value = 'PASS'
$ws.Cells.Item(115, 5).Value = '[]'
$ws.Cells.Item(116, 1).Value = '(''FR54 6321 1036 07CL FIPI W3AJ A72'', '''')'
$ws.Cells.Item(116, 2).Value = '[2015.0, 31.6075, 22.79953, 24.24]'
$ws.Cells.Item(116, 3).Value = '[2015.0, 31.6075, 22.79953, 24.24]'
$ws.Cells.Item(116, 4).Value = 'PASS'
$ws.Cells.Item(116, 5).Value = '[]'
$ws.Cells.Item(117, 1).Value = '(''GL16 3614 7192 7269 45'', '''')'
$ws.Cells.Item(117, 2).Value = '[2018.0, 6.71129, 5.96627, 15.6]'
$ws.Cells.Item(117, 3).Value = '[2018.0, 6.71129, 5.96627, 15.6]'
$ws.Cells.Item(117, 4).Value = 'PASS'
$ws.Cells.Item(117, 5).Value = '[]'
$ws.Cells.Item(118, 1).Value = '(''CH31 5106 86EY EGNR WZ9U L'', '''')'
$ws.Cells.Item(118, 2).Value = '[2019.0, 17.51596, 16.62937, 17.69]'
$ws.Cells.Item(118, 3).Value = '[2019.0, 17.51596, 16.62937, 17.69]'
$ws.Cells.Item(118, 4).Value = 'PASS'
$ws.Cells.Item(118, 5).Value = '[]'
$ws.Cells.Item(119, 1).Value = '(''GR30 9103 924T 492V ODJA MC2H V73'', ''matrix'')'
$ws.Cells.Item(119, 2).Value = '[2016.0, 21.39719, 9.31935, 8.07]'
$ws.Cells.Item(119, 3).Value = '[2016.0, 21.39719, 9.31935, 8.07]'
$ws.Cells.Item(119, 4).Value = 'PASS'
$ws.Cells.Item(119, 5).Value = '[]'
$ws.Cells.Item(120, 1).Value = '(''EE60 0860 2937 4213 4766'', ''Enhanced'')'
$ws.Cells.Item(120, 2).Value = '[2013.0, 33.92035, 17.96162, 13.88]'
$ws.Cells.Item(120, 3).Value = '[2013.0, 33.92035, 17.96162, 13.88]'
$ws.Cells.Item(120, 4).Value = 'PASS'
$ws.Cells.Item(120, 5).Value = '[]'
$ws.Cells.Item(121, 1).Value = '(''PT48 9547 7145 0711 6406 9763 0'', ''Quality-focused'')'
$ws.Cells.Item(121, 2).Value = '[2020.0, 25.34973, 30.41719, 9.06]'
$ws.Cells.Item(121, 3).Value = '[2020.0, 25.34973, 30.41719, 9.06]'
$ws.Cells.Item(121, 4).Value = 'PASS'
$ws.Cells.Item(121, 5).Value = '[]'
$ws.Cells.Item(122, 1).Value = '(''FR44 6285 3462 75FF SDKU 8JEE Z98'', '''')'
$ws.Cells.Item(122, 2).Value = '[2015.0, 15.86295, 30.64175, 13.3]'
$ws.Cells.Item(122, 3).Value = '[2015.0, 15.86295, 30.64175, 13.3]'
$ws.Cells.Item(122, 4).Value = 'PASS'
$ws.Cells.Item(122, 5).Value = '[]'
$ws.Cells.Item(123, 1).Value = '(''FR16 0580 2045 20MG KTEV ZQNO 396'', '''')'
$ws.Cells.Item(123, 2).Value = '[2010.0, 19.01096, 18.36853, 32.16]'
$ws.Cells.Item(123, 3).Value = '[2010.0, 19.01096, 18.36853, 32.16]'
$ws.Cells.Item(123, 4).Value = 'PASS'
$ws.Cells.Item(123, 5).Value = '[]'
$ws.Cells.Item(124, 1).Value = '(''AD11 8389 0233 NS7R WURI 9OPT'', ''Progressive'')'
$ws.Cells.Item(124, 2).Value = '[2004.0, 8.04103, 4.39201, 2.61]'
$ws.Cells.Item(124, 3).Value = '[2004.0, 8.04103, 4.39201, 2.61]'
$ws.Cells.Item(124, 4).Value = 'PASS'
$ws.Cells.Item(124, 5).Value = '[]'
$ws.Cells.Item(125, 1).Value = '(''SE73 9240 7477 7911 9057 9462'', '''')'
$ws.Cells.Item(125, 2).Value = '[2008.0, 12.53634, 28.19933, 34.03]'
$ws.Cells.Item(125, 3).Value = '[2008.0, 12.53634, 28.19933, 34.03]'
$ws.Cells.Item(125, 4).Value = 'PASS'
$ws.Cells.Item(125, 5).Value = '[]'
$ws.Cells.Item(126, 1).Value = '(''PS17 YVSS UAFQ VR9I STIU TL70 PVBY I'', '''')'
$ws.Cells.Item(126, 2).Value = '[2011.0, 3.40426, 14.73376, 26.67]'
$ws.Cells.Item(126, 3).Value = '[2011.0, 3.40426, 14.73376, 26.67]'
$ws.Cells.Item(126, 4).Value = 'PASS'
$ws.Cells.Item(126, 5).Value = '[]'
$ws.Cells.Item(127, 1).Value = '(''MD50 ZOFQ MPTQ MV5P KHHR 8TBN'', ''User-friendly'')'
$ws.Cells.Item(127, 2).Value = '[2014.0, 4.43078, 20.29182, 30.96]'
$ws.Cells.Item(127, 3).Value = '[2014.0, 4.43078, 20.29182, 30.96]'
$ws.Cells.Item(127, 4).Value = 'PASS'
$ws.Cells.Item(127, 5).Value = '[]'
$ws.Cells.Item(128, 1).Value = '(''FI49 6049 5956 9898 35'', ''neutral'')'
$ws.Cells.Item(128, 2).Value = '[2017.0, 25.06377, 33.71366, 22.04]'
$ws.Cells.Item(128, 3).Value = '[2017.0, 25.06377, 33.71366, 22.04]'
$ws.Cells.Item(128, 4).Value = 'PASS'
$ws.Cells.Item(128, 5).Value = '[]'
$ws.Cells.Item(129, 1).Value = '(''FR26 4582 5102 64DY 4DY6 LI7G E21'', ''Automated'')'
$ws.Cells.Item(129, 2).Value = '[2007.0, 20.82924, 14.15138, 16.84]'
$ws.Cells.Item(129, 3).Value = '[2007.0, 20.82924, 14.15138, 16.84]'
$ws.Cells.Item(129, 4).Value = 'PASS'
$ws.Cells.Item(129, 5).Value = '[]'
$ws.Cells.Item(130, 1).Value = '(''LI72 7685 1IGO GJCT HSAJ K'', ''Graphical User Interface'')'
$ws.Cells.Item(130, 2).Value = '[2013.0, 18.02644, 7.84381, 1.09]'
$ws.Cells.Item(130, 3).Value = '[2013.0, 18.02644, 7.84381, 1.09]'
$ws.Cells.Item(130, 4).Value = 'PASS'
$ws.Cells.Item(130, 5).Value = '[]'
$ws.Cells.Item(131, 1).Value = '(''AD39 7446 9930 MWO3 CIXZ JXUM'', ''middleware'')'
$ws.Cells.Item(131, 2).Value = '[2011.0, 33.74886, 17.95307, 2.94]'
$ws.Cells.Item(131, 3).Value = '[2011.0, 33.74886, 17.95307, 2.94]'
$ws.Cells.Item(131, 4).Value = 'PASS'
$ws.Cells.Item(131, 5).Value = '[]'
$ws.Cells.Item(132, 1).Value = '(''PT86 9472 7150 6598 9677 1790 8'', '''')'
$ws.Cells.Item(132, 2).Value = '[2004.0, 11.99216, 33.11774, 10.6]'
$ws.Cells.Item(132, 3).Value = '[2004.0, 11.99216, 33.11774, 10.6]'
$ws.Cells.Item(132, 4).Value = 'PASS'
$ws.Cells.Item(132, 5).Value = '[]'
$ws.Cells.Item(133, 1).Value = '(''CY92 9856 5900 JNUI MQDJ GMZ5 SXY3'', '''')'
$ws.Cells.Item(133, 2).Value = '[2013.0, 14.08251, 11.45311, 5.92]'
$ws.Cells.Item(133, 3).Value = '[2013.0, 14.08251, 11.45311, 5.92]'
$ws.Cells.Item(133, 4).Value = 'PASS'
$ws.Cells.Item(133, 5).Value = '[]'
$ws.Cells.Item(134, 1).Value = '(''MT05 TBDE 5852 9M5Y OEAO MA4U LBHG A9N'', '''')'
$ws.Cells.Item(134, 2).Value = '[2011.0, 23.46623, 27.75647, 21.45]'
$ws.Cells.Item(134, 3).Value = '[2011.0, 23.46623, 27.75647, 21.45]'
$ws.Cells.Item(134, 4).Value = 'PASS'
$ws.Cells.Item(134, 5).Value = '[]'
$ws.Cells.Item(135, 1).Value = '(''PS22 XCSM A4JV ZWHN GKJH NN0I SJAV N'', ''Synergistic'')'
$ws.Cells.Item(135, 2).Value = '[2004.0, 6.60445, 32.82567, 14.59]'
$ws.Cells.Item(135, 3).Value = '[2004.0, 6.60445, 32.82567, 14.59]'
$ws.Cells.Item(135, 4).Value = 'PASS'
$ws.Cells.Item(135, 5).Value = '[]'
$ws.Cells.Item(136, 1).Value = '(''MK40 019S AXFJ AHYJ D76'', '''')'
$ws.Cells.Item(136, 2).Value = '[2011.0, 31.03655, 7.06158, 1.82]'
$ws.Cells.Item(136, 3).Value = '[2011.0, 31.03655, 7.06158, 1.82]'
$ws.Cells.Item(136, 4).Value = 'PASS'
$ws.Cells.Item(136, 5).Value = '[]'
$ws.Cells.Item(137, 1).Value = '(''LU85 097J NGB4 CYVN 067E'', '''')'
$ws.Cells.Item(137, 2).Value = '[2012.0, 18.80157, 17.98329, 9.32]'
$ws.Cells.Item(137, 3).Value = '[2012.0, 18.80157, 17.98329, 9.32]'
$ws.Cells.Item(137, 4).Value = 'PASS'
$ws.Cells.Item(137, 5).Value = '[]'
$ws.Cells.Item(138, 1).Value = '(''PK88 MODL G7O1 HQN9 1GPW X3XN'', '''')'
$ws.Cells.Item(138, 2).Value = '[2016.0, 24.1303, 19.53632, 4.3]'
$ws.Cells.Item(138, 3).Value = '[2016.0, 24.1303, 19.53632, 4.3]'
$ws.Cells.Item(138, 4).Value = 'PASS'
$ws.Cells.Item(138, 5).Value = '[]'
$ws.Cells.Item(139, 1).Value = '(''ES22 9700 7734 8770 4091 5098'', '''')'
$ws.Cells.Item(139, 2).Value = '[2009.0, 6.07865, 24.7863, 5.22]'
$ws.Cells.Item(139, 3).Value = '[2009.0, 6.07865, 24.7863, 5.22]'
$ws.Cells.Item(139, 4).Value = 'PASS'
$ws.Cells.Item(139, 5).Value = '[]'
$ws.Cells.Item(140, 1).Value = '(''AT93 5616 4483 1930 1613'', '''')'
$ws.Cells.Item(140, 2).Value = '[2017.0, 25.15557, 17.36863, 10.99]'
$ws.Cells.Item(140, 3).Value = '[2017.0, 25.15557, 17.36863, 10.99]'
$ws.Cells.Item(140, 4).Value = 'PASS'
$ws.Cells.Item(140, 5).Value = '[]'
$ws.Cells.Item(141, 1).Value = '(''FR54 4415 3882 12MY GWIB WULK D13'', '''')'
$ws.Cells.Item(141, 2).Value = '[2013.0, 23.52919, 14.31729, 26.44]'
$ws.Cells.Item(141, 3).Value = '[2013.0, 23.52919, 14.31729, 26.44]'
$ws.Cells.Item(141, 4).Value = 'PASS'
$ws.Cells.Item(141, 5).Value = '[]'
$ws.Cells.Item(142, 1).Value = '(''HU85 1357 8508 2949 0123 4827 6139'', '''')'
$ws.Cells.Item(142, 2).Value = '[2013.0, 21.2783, 33.1919, 29.54]'
$ws.Cells.Item(142, 3).Value = '[2013.0, 21.2783, 33.1919, 29.54]'
$ws.Cells.Item(142, 4).Value = 'PASS'
$ws.Cells.Item(142, 5).Value = '[]'
$ws.Cells.Item(143, 1).Value = '(''LV91 WMUV JOFO 2XLR UYEI T'', ''Innovative'')'
$ws.Cells.Item(143, 2).Value = '[2008.0, 28.97925, 6.71373, 29.25]'
$ws.Cells.Item(143, 3).Value = '[2008.0, 28.97925, 6.71373, 29.25]'
$ws.Cells.Item(143, 4).Value = 'PASS'
$ws.Cells.Item(143, 5).Value = '[]'
$ws.Cells.Item(144, 1).Value = '(''FR06 2928 8076 28DQ 2LU7 4FID F04'', ''didactic'')'
$ws.Cells.Item(144, 2).Value = '[2002.0, 24.67853, 25.10648, 23.27]'
$ws.Cells.Item(144, 3).Value = '[2002.0, 24.67853, 25.10648, 23.27]'
$ws.Cells.Item(144, 4).Value = 'PASS'
$ws.Cells.Item(144, 5).Value = '[]'
$ws.Cells.Item(145, 1).Value = '(''AL97 9974 2251 MJPR C5VC VVBS VZMC'', '''')'
$ws.Cells.Item(145, 2).Value = '[2012.0, 34.75131, 5.40346, 12.32]'
$ws.Cells.Item(145, 3).Value = '[2012.0, 34.75131, 5.40346, 12.32]'
$ws.Cells.Item(145, 4).Value = 'PASS'
$ws.Cells.Item(145, 5).Value = '[]'
$ws.Cells.Item(146, 1).Value = '(''HU52 9898 4213 5176 1777 7165 5419'', ''encompassing'')'
$ws.Cells.Item(146, 2).Value = '[2012.0, 1.7833, 23.57441, 23.3]'
$ws.Cells.Item(146, 3).Value = '[2012.0, 0.7833, 23.57441, 23.3]'
$ws.Cells.Item(146, 4).Value = 'FAIL'
$ws.Cells.Item(146, 5).Value = '[-1.0]'
$ws.Cells.Item(147, 1).Value = '(''AT12 1939 6455 5066 7889'', '''')'
$ws.Cells.Item(147, 2).Value = '[2014.0, 2.57696, 27.78488, 13.99]'
$ws.Cells.Item(147, 3).Value = '[2014.0, 2.57696, 27.78488, 13.99]'
$ws.Cells.Item(147, 4).Value = 'PASS'
$ws.Cells.Item(147, 5).Value = '[]'
$ws.Cells.Item(148, 1).Value = '(''SK05 5786 1420 0506 2569 0799'', '''')'
$ws.Cells.Item(148, 2).Value = '[2006.0, 5.37354, 15.04885, 16.5]'
$ws.Cells.Item(148, 3).Value = '[2006.0, 5.37354, 15.04885, 16.5]'
$ws.Cells.Item(148, 4).Value = 'PASS'
$ws.Cells.Item(148, 5).Value = '[]'
$ws.Cells.Item(149, 1).Value = '(''MC25 6140 6612 40NT JOOH KYXP H72'', ''local'')'
$ws.Cells.Item(149, 2).Value = '[2000.0, 9.85837, 6.05067, 31.17]'
$ws.Cells.Item(149, 3).Value = '[2000.0, 9.85837, 6.05067, 31.17]'
$ws.Cells.Item(149, 4).Value = 'PASS'
$ws.Cells.Item(149, 5).Value = '[]'

# Fix conditional PASS/FAIL cell styles (fill color) to follow the moved rows.
# Use a source cell that currently has FAIL formatting (red fill) to copy format from,
# and a source cell that currently has PASS formatting (green fill) likewise.
$failFormatSource = $ws.Range("D27")
$passFormatSource = $ws.Range("D2")

$failFormatSource.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D75").PasteSpecial(-4122)
$ws.Range("D146").PasteSpecial(-4122)

$passFormatSource.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("D112").PasteSpecial(-4122)
$ws.Range("D141").PasteSpecial(-4122)
